# "Generate Report for Archive"
#
# The localization-status report is regenerated: the three still-pending
# files (10707ff1, 5a306eb0, ac4d99b3) are re-sorted in every sheet so the
# file that is furthest along (10707ff1, "Ready for handoff") sinks to the
# bottom of the pending block, while the two files that have since moved
# into translation (5a306eb0, ac4d99b3) bubble up and their status flips
# from "Ready for handoff" to "In Translation".
#
# Concretely, for rows 3/4/5 of every sheet:
#   new row3 <- old row4 data, status -> "In Translation"
#   new row4 <- old row5 data, status -> "In Translation"
#   new row5 <- old row3 data, status unchanged ("Ready for handoff")

$wb = $excel.ActiveWorkbook

function Get-RowValues($ws, [int]$row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range($c + $row).Value2
    }
    return $vals
}

function Set-RowValues($ws, [int]$row, $cols, $vals) {
    foreach ($c in $cols) {
        $v = $vals[$c]
        if ($v -ne "") {
            $ws.Range($c + $row).Value = $v
        }
    }
}

function Set-HyperlinkDisplay($ws, [int]$row, [string]$col, [string]$newText) {
    foreach ($hl in $ws.Hyperlinks) {
        if (($hl.Range.Row -eq $row) -and ($hl.Range.Column -eq $col)) {
            $hl.TextToDisplay = $newText
        }
    }
}

# ---------- Sheet "Overview" ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewCols = @("A","B","C","D")

$ovRow3 = Get-RowValues $wsOverview 3 $overviewCols
$ovRow4 = Get-RowValues $wsOverview 4 $overviewCols
$ovRow5 = Get-RowValues $wsOverview 5 $overviewCols

$newOvRow3 = $ovRow4.Clone()
$newOvRow3["B"] = "In Translation"
$newOvRow3["C"] = "In Translation"

$newOvRow4 = $ovRow5.Clone()
$newOvRow4["B"] = "In Translation"
$newOvRow4["C"] = "In Translation"

$newOvRow5 = $ovRow3.Clone()

Set-RowValues $wsOverview 3 $overviewCols $newOvRow3
Set-RowValues $wsOverview 4 $overviewCols $newOvRow4
Set-RowValues $wsOverview 5 $overviewCols $newOvRow5

Set-HyperlinkDisplay $wsOverview 3 1 $newOvRow3["A"]
Set-HyperlinkDisplay $wsOverview 4 1 $newOvRow4["A"]
Set-HyperlinkDisplay $wsOverview 5 1 $newOvRow5["A"]

# ---------- Per-language sheets "zh-cn" / "de-de" ----------
$langSheets = @("zh-cn", "de-de")
$langCols = @("A","B","C","D","E","F","G","H","I","J","K","L")

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row3 = Get-RowValues $ws 3 $langCols
    $row4 = Get-RowValues $ws 4 $langCols
    $row5 = Get-RowValues $ws 5 $langCols

    $newRow3 = $row4.Clone()
    $newRow3["C"] = "In Translation"

    $newRow4 = $row5.Clone()
    $newRow4["C"] = "In Translation"

    $newRow5 = $row3.Clone()

    Set-RowValues $ws 3 $langCols $newRow3
    Set-RowValues $ws 4 $langCols $newRow4
    Set-RowValues $ws 5 $langCols $newRow5

    Set-HyperlinkDisplay $ws 3 1 $newRow3["A"]
    Set-HyperlinkDisplay $ws 3 4 $newRow3["D"]
    Set-HyperlinkDisplay $ws 4 1 $newRow4["A"]
    Set-HyperlinkDisplay $ws 4 4 $newRow4["D"]
    Set-HyperlinkDisplay $ws 5 1 $newRow5["A"]
    Set-HyperlinkDisplay $ws 5 4 $newRow5["D"]
}
